$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    102 = @('HSBCDOL', 'BNF', 'Acciones de Sociedades de Inversion', 'Mercado de Deuda Extranjero', 'Renta Fija Internacional')
    103 = @('NTEDLS+', 'FF', 'Acciones de Sociedades de Inversion', 'Mercado de Deuda Extranjero', 'Renta Fija Internacional')
    104 = @('PRGLOB', 'FFR', 'Acciones de Sociedades de Inversion', 'Mercado de Deuda Extranjero', 'Renta Fija Internacional')
    105 = @('SCOTDL+', 'C1E', 'Acciones de Sociedades de Inversion', 'Mercado de Deuda Extranjero', 'Renta Fija Internacional')
    106 = @('SCOTDOL', 'C1E', 'Acciones de Sociedades de Inversion', 'Mercado de Deuda Extranjero', 'Renta Fija Internacional')
    107 = @('SURUSD', 'BOE', 'Acciones de Sociedades de Inversion', 'Mercado de Deuda Extranjero', 'Renta Fija Internacional')
    108 = @('PRINHYD', 'FFX', 'Acciones de Sociedades de Inversion de Instrumentos de Deuda', 'Mercado de Deuda Extranjero', 'Renta Fija Internacional')
    109 = @('NTEINT+', 'FF', 'Acciones de Sociedades de Inversion', 'Mercado de Capitales Extranjero', 'Renta Variable Internacional')
    110 = @('NTEUSA', 'FF', 'Acciones de Sociedades de Inversion', 'Mercado de Capitales Extranjero', 'Renta Variable Internacional')
    111 = @('NTEUSA+', 'FF', 'Acciones de Sociedades de Inversion', 'Mercado de Capitales Extranjero', 'Renta Variable Internacional')
    112 = @('SCOTGLO', 'C1E', 'Acciones de Sociedades de Inversion', 'Mercado de Capitales Extranjero', 'Renta Variable Internacional')
    113 = @('SURGLOB', 'BOE0', 'Acciones de Sociedades de Inversion', 'Mercado de Capitales Extranjero', 'Renta Variable Internacional')
    114 = @('BLKINT1', 'M0-A', 'Acciones de Sociedades de Inversion', 'Mercado de Capitales Extranjero', 'Renta Variable Internacional')
    115 = @('FT-GLOB', 'BE', 'Acciones de Sociedades de Inversion', 'Mercado de Capitales Extranjero', 'Renta Variable Internacional')
    116 = @('PEMERGE', 'FFR', 'Acciones de Sociedades de Inversion', 'Mercado de Capitales Extranjero', 'Renta Variable Internacional')
    117 = @('PRINFUS', 'FFX', 'Acciones de Sociedades de Inversion', 'Mercado de Capitales Extranjero', 'Renta Variable Internacional')
}

foreach ($r in ($data.Keys | Sort-Object {[int]$_})) {
    $rowVals = $data[$r]
    for ($c = 1; $c -le 5; $c++) {
        $ws.Cells.Item([int]$r, $c).Value = $rowVals[$c - 1]
    }
    $ws.Range($ws.Cells.Item([int]$r, 1), $ws.Cells.Item([int]$r, 5)).Interior.Color = 65535
}

$ws.Range("B105").Select()
